$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 249
$ws.Range("I29").Value = 249
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 747
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -466
$ws.Range("N29").Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3665
$ws.Range("I106").Value = 3665
$ws.Range("K106").Value = 3665
$ws.Range("M106").Value = -3034

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14094.348
$ws.Range("I132").Value = 2764.5
$ws.Range("K132").Value = 8293.5
$ws.Range("M132").Value = -5763.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2808.0303
$ws.Range("I138").Value = 1471.6923
$ws.Range("J138").Value = 3676.65
$ws.Range("K138").Value = 4415.0769
$ws.Range("L138").Value = 11029.95
$ws.Range("M138").Value = 724.9231
$ws.Range("N138").Value = -21309.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 218
$ws.Range("I5").Value = 218
$ws.Range("K5").Value = 218
$ws.Range("M5").Value = -106

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15154387
$ws.Range("I32").Value = 15875857
$ws.Range("K32").Value = 15875857
$ws.Range("M32").Value = -15875570

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2734.25
$ws.Range("I45").Value = 1602
$ws.Range("K45").Value = 1602
$ws.Range("M45").Value = -1225

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1973.4
$ws.Range("I61").Value = 2482.2
$ws.Range("J61").Value = 955.8
$ws.Range("K61").Value = 2482.2
$ws.Range("L61").Value = 955.8
$ws.Range("M61").Value = -2270.2
$ws.Range("N61").Value = -1379.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1551.1666
$ws.Range("I74").Value = 1669.0303
$ws.Range("K74").Value = 1669.0303
$ws.Range("M74").Value = -795.0302999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1551.1666
$ws.Range("I77").Value = 1669.0303
$ws.Range("K77").Value = 8345.1515
$ws.Range("M77").Value = -3977.1515

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4183.5
$ws.Range("I102").Value = 4076.647
$ws.Range("K102").Value = 4076.647
$ws.Range("M102").Value = -2454.647

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3717.4849
$ws.Range("I122").Value = 2876.8696
$ws.Range("K122").Value = 8630.6088
$ws.Range("M122").Value = -6180.6088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1973.4
$ws.Range("I136").Value = 2482.2
$ws.Range("J136").Value = 955.8
$ws.Range("K136").Value = 7446.599999999999
$ws.Range("L136").Value = 2867.4
$ws.Range("M136").Value = -4896.599999999999
$ws.Range("N136").Value = -7967.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 218
$ws.Range("I4").Value = 218
$ws.Range("K4").Value = 218
$ws.Range("M4").Value = -103

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2004.1305
$ws.Range("I20").Value = 1322.4
$ws.Range("K20").Value = 1322.4
$ws.Range("M20").Value = -1075.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1892.75
$ws.Range("I105").Value = 1840.2667
$ws.Range("J105").Value = 2050.2
$ws.Range("K105").Value = 1840.2667
$ws.Range("L105").Value = 2050.2
$ws.Range("M105").Value = -93.2666999999999
$ws.Range("N105").Value = -5544.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2410.3044
$ws.Range("I134").Value = 2164.611
$ws.Range("J134").Value = 3294.8
$ws.Range("K134").Value = 6493.833
$ws.Range("L134").Value = 9884.400000000001
$ws.Range("M134").Value = -3958.833
$ws.Range("N134").Value = -14954.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 122394.89
$ws.Range("I4").Value = 122394.89
$ws.Range("K4").Value = 122394.89
$ws.Range("M4").Value = -122282.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2046.1
$ws.Range("I58").Value = 1401.375
$ws.Range("J58").Value = 4625
$ws.Range("K58").Value = 1401.375
$ws.Range("L58").Value = 4625
$ws.Range("M58").Value = -1198.375
$ws.Range("N58").Value = -5031

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2523.2778
$ws.Range("I132").Value = 2078.1
$ws.Range("J132").Value = 4749.1665
$ws.Range("K132").Value = 6234.299999999999
$ws.Range("L132").Value = 14247.4995
$ws.Range("M132").Value = -3704.299999999999
$ws.Range("N132").Value = -19307.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10999.125
$ws.Range("I134").Value = 10282.714
$ws.Range("K134").Value = 30848.142
$ws.Range("M134").Value = -28313.142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2046.1
$ws.Range("I136").Value = 1401.375
$ws.Range("J136").Value = 4625
$ws.Range("K136").Value = 4204.125
$ws.Range("L136").Value = 13875
$ws.Range("M136").Value = -1654.125
$ws.Range("N136").Value = -18975

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 61078950
$ws.Range("I4").Value = 62987664
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 188962992
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = -188962880
$ws.Range("N4").Value = -974

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2420.3333
$ws.Range("I94").Value = 297.33334
$ws.Range("J94").Value = 6666.3335
$ws.Range("K94").Value = 892.0000200000001
$ws.Range("L94").Value = 19999.0005
$ws.Range("M94").Value = -216.0000200000001
$ws.Range("N94").Value = -21351.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2194
$ws.Range("I121").Value = 916.6667
$ws.Range("J121").Value = 2741.4285
$ws.Range("K121").Value = 2750.0001
$ws.Range("L121").Value = 8224.2855
$ws.Range("M121").Value = -1440.0001
$ws.Range("N121").Value = -10844.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4138.2144
$ws.Range("I131").Value = 3025
$ws.Range("J131").Value = 4583.5
$ws.Range("K131").Value = 9075
$ws.Range("L131").Value = 13750.5
$ws.Range("M131").Value = -4035
$ws.Range("N131").Value = -23830.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1975.9375
$ws.Range("J137").Value = 3227.4
$ws.Range("L137").Value = 9682.200000000001
$ws.Range("N137").Value = -19882.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7622.2
$ws.Range("I122").Value = 7244.6
$ws.Range("K122").Value = 21733.8
$ws.Range("M122").Value = -19283.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5673.0586
$ws.Range("I132").Value = 5808.875
$ws.Range("K132").Value = 17426.625
$ws.Range("M132").Value = -14896.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8700.058999999999
$ws.Range("I40").Value = 10639
$ws.Range("K40").Value = 10639
$ws.Range("M40").Value = -10503

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 386.6875
$ws.Range("I55").Value = 583
$ws.Range("J55").Value = 190.375
$ws.Range("K55").Value = 583
$ws.Range("L55").Value = 190.375
$ws.Range("M55").Value = -410
$ws.Range("N55").Value = -536.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3078.5144
$ws.Range("I93").Value = 2082.96
$ws.Range("J93").Value = 5567.4
$ws.Range("K93").Value = 2082.96
$ws.Range("L93").Value = 5567.4
$ws.Range("M93").Value = -834.96
$ws.Range("N93").Value = -8063.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 270992.25
$ws.Range("J116").Value = 270992.25
$ws.Range("L116").Value = 270992.25
$ws.Range("N116").Value = -280170.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 64477.273
$ws.Range("J127").Value = 70000
$ws.Range("L127").Value = 70000
$ws.Range("N127").Value = -79920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2492.0476
$ws.Range("J132").Value = 3530.7273
$ws.Range("L132").Value = 10592.1819
$ws.Range("N132").Value = -15652.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19788.8
$ws.Range("I2").Value = 20986.25
$ws.Range("J2").Value = 14999
$ws.Range("K2").Value = 20986.25
$ws.Range("L2").Value = 14999
$ws.Range("M2").Value = -20874.25
$ws.Range("N2").Value = -15223

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 22124.75
$ws.Range("I52").Value = 4250
$ws.Range("J52").Value = 39999.5
$ws.Range("K52").Value = 4250
$ws.Range("L52").Value = 39999.5
$ws.Range("M52").Value = -4024
$ws.Range("N52").Value = -40451.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 936.1429000000001
$ws.Range("I107").Value = 897.4286
$ws.Range("J107").Value = 974.8570999999999
$ws.Range("K107").Value = 2692.2858
$ws.Range("L107").Value = 2924.5713
$ws.Range("M107").Value = -772.2857999999997
$ws.Range("N107").Value = -6764.5713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 632.0909
$ws.Range("I113").Value = 583.7222
$ws.Range("J113").Value = 849.75
$ws.Range("K113").Value = 1751.1666
$ws.Range("L113").Value = 2549.25
$ws.Range("M113").Value = 418.8334
$ws.Range("N113").Value = -6889.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1805.8182
$ws.Range("I122").Value = 1681.3529
$ws.Range("J122").Value = 2229
$ws.Range("K122").Value = 5044.0587
$ws.Range("L122").Value = 6687
$ws.Range("M122").Value = -2594.0587
$ws.Range("N122").Value = -11587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 914.91895
$ws.Range("I132").Value = 914.91895
$ws.Range("K132").Value = 2744.75685
$ws.Range("M132").Value = -214.7568499999998
